$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number that was updated
# from 45190 (2023-09-21) to 45192 (2023-09-23) for every data row
# (rows 2 through 360).
$ws.Range("C2:C360").Value = (Get-Date -Year 2023 -Month 9 -Day 23 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)
